$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Replace the long INSABA source citation strings in A31 and A33
# with the short "INSABA" label (matching A30 / A32).
$ws.Range("A31").Value = "INSABA"
$ws.Range("A33").Value = "INSABA"
